$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rng = $ws.Range("A30:F31")
$rng.NumberFormat = "@"

# Row 30
$ws.Cells.Item(30, 1).Value = "20140720"
$ws.Cells.Item(30, 2).Value = "CMPN OCHEESECAKE 90"
$ws.Cells.Item(30, 3).Value = "RCI04N"
$ws.Cells.Item(30, 4).Value = "2"
$ws.Cells.Item(30, 5).Value = "12"
$ws.Cells.Item(30, 6).Value = "RT,(E-1B)"

# Row 31
$ws.Cells.Item(31, 1).Value = "20140719"
$ws.Cells.Item(31, 2).Value = "CMPNA BON VAN CHO 45"
$ws.Cells.Item(31, 3).Value = "RCI04N"
$ws.Cells.Item(31, 4).Value = "2"
$ws.Cells.Item(31, 5).Value = "13"
$ws.Cells.Item(31, 6).Value = "RT,(E-1B)"

$ws.Range("A30:F31").Borders.LineStyle = 1
